$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - update metric values
$ws.Range("B3").Value = 0.9964225133488052
$ws.Range("C3").Value = 0.9961962861686517
$ws.Range("D3").Value = 0.9951104830864469

# Row 4: rename model and update metric values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9956952833923552
$ws.Range("C4").Value = 0.9955967098746941
$ws.Range("D4").Value = 0.9771262908096249

# Row 5: rename model and update metric values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9974406134358692
$ws.Range("C5").Value = 0.9968198735030032
$ws.Range("D5").Value = 0.9963345391748396
